$d = $word.ActiveDocument

# Locate the exact text "kod" that needs to be split/extended.
$rng = $d.Content
$found = $rng.Find.Execute("kod", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' +
           '<w:r><w:t>K</w:t></w:r>' +
           '<w:r><w:t>od</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> 2. versiyon</w:t></w:r>' +
           '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}
